$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting rows 30-45 down to 31-46
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with fresh data
$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(30, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(30, 4).Value = 44455
$ws.Cells.Item(30, 4).NumberFormat = $ws.Cells.Item(31, 4).NumberFormat
$ws.Cells.Item(30, 5).Value = 15
$ws.Cells.Item(30, 6).Value = 100112038
$ws.Cells.Item(30, 7).Value = "Cebollín baby"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 300
$ws.Cells.Item(30, 11).Value = 2000
$ws.Cells.Item(30, 12).Value = 2500
$ws.Cells.Item(30, 13).Value = 2250
$ws.Cells.Item(30, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(30, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(30, 16).Value = 1125
$ws.Cells.Item(30, 17).Value = 2
$ws.Cells.Item(30, 18).Value = "Hortaliza"
